$d = $word.ActiveDocument

# 1) Merge "IF a>=75" onto the line above by removing the " 0<=a<=100" line
#    (together with its trailing line break) and adding " Then" to the
#    condition, so the paragraph reads "...INPUT a<br>IF a>=75 Then<br>...".
$d.Content.Find.Execute(
    " 0<=a<=100" + [char]11 + "IF a>=75",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "IF a>=75 Then", 2) | Out-Null

# 2) "IF 65<=a<75 " -> "IF 65<=a Then"
$d.Content.Find.Execute(
    "IF 65<=a<75 ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "IF 65<=a Then", 2) | Out-Null

# 3) "IF 45<=a<65" -> "IF 45<=a Then"
$d.Content.Find.Execute(
    "IF 45<=a<65",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "IF 45<=a Then", 2) | Out-Null

# 4) "IF 35<=a<45" -> "IF 35<=a Then". A placeholder character ("@") is kept
#    right after "Then" for now - it gives us a safe (non paragraph-end)
#    anchor to drop the moved "_GoBack" bookmark on, since collapsed
#    ranges exactly at a paragraph boundary aren't placed reliably. The
#    placeholder is stripped again immediately afterwards.
$d.Content.Find.Execute(
    "IF 35<=a<45",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "IF 35<=a Then@", 2) | Out-Null

$ifPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "IF 35<=a Then@") {
        $ifPara = $p
    }
}
$bmPos = $ifPara.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$d.Content.Find.Execute(
    "@",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 2) | Out-Null

# 5) "Display D" -> "Display Loai D"
$d.Content.Find.Execute(
    "Display D",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Display Loai D", 2) | Out-Null

# 6) "Else Display E" -> "Else Display Loai E"
$d.Content.Find.Execute(
    "Else Display E",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Else Display Loai E", 2) | Out-Null
